# Add a new "UK" Test-Data sheet, cloned from the last existing country
# sheet ("Poland"), with UK-specific values - mirrors how every other
# per-country sheet in this workbook was produced.

$wb = $excel.ActiveWorkbook

# Clone the last country sheet ("Poland") and drop the copy at the very
# end of the tab strip (after Poland), then rename it to "UK".
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$template = $wb.Worksheets.Item("Poland")
$template.Copy($null, $lastSheet)

$uk = $wb.Worksheets.Item($wb.Worksheets.Count)
$uk.Name = "UK"

# The UK sheet has one extra "system" row (GMPIM) that the other markets
# don't have, inserted right above the PR1D2/Wg/Miscellaneous rows - this
# pushes everything below it down by one row.
$uk.Rows.Item(9).Insert()

# Bring the formatting for the newly-inserted row in line with its
# neighbours (A8/A10 both use the boxed "code list" style).
$uk.Range("A8").Copy()
$uk.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the UK-specific content.
$uk.Range("A9").Value = "GMPIM"
$uk.Range("B4").Value = "NGC-2741/T3343/T3342/T3345"
$uk.Range("B2").Value = "UK Market"

# Make the new sheet the active one, with A9 selected - matching how it
# was left after entering the data.
$uk.Activate()
$uk.Range("A9").Select() | Out-Null
